$wb = $excel.ActiveWorkbook

# --- Datos_Limpios: change the product sold in row 4 from "Escritorio" to "Laptop" ---
$wsDatos = $wb.Worksheets.Item("Datos_Limpios")
$wsDatos.Range("H4").Value = "Laptop"
$wsDatos.Range("H4").Select() | Out-Null

# --- Analisis: add a new "Ventas de Tecnologia en Bogota" block below the existing tables ---
$wsAnalisis = $wb.Worksheets.Item("Analisis")
$wsAnalisis.Activate() | Out-Null

$wsAnalisis.Range("A8").Value = "Ventas de Tecnologia en Bogota"
$wsAnalisis.Range("A8").WrapText = $true
$wsAnalisis.Range("A8").VerticalAlignment = -4108
$wsAnalisis.Rows.Item(8).RowHeight = 30

$wsAnalisis.Range("A9").Formula = '=SUMIFS(tbl_Ventas[Total_Ventas],tbl_Ventas[Categoria],"Tecnología",tbl_Ventas[Ciudad_Cliente],"Bogotá")'

$wsAnalisis.Range("A10").Select() | Out-Null
